# Weekly update: a new price record is inserted right after the existing
# row 233 (i.e. at row 234), pushing all subsequent "Jengibre" rows down by
# one (old row 234 -> new row 235, ..., old row 337 -> new row 338).
# The sheet's used range grows from A1:R337 to A1:R338 automatically as a
# result of the row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 234; everything from the old row 234 onward
# shifts down one row (old 234 -> 235, ..., old 337 -> 338).
$ws.Rows("234").Insert()

# Populate the newly inserted row 234 with this week's record.
$ws.Cells.Item(234, 1).Value  = 10
$ws.Cells.Item(234, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(234, 3).Value  = "La Araucanía"
$ws.Cells.Item(234, 4).Value  = 45141
$ws.Cells.Item(234, 5).Value  = 9
$ws.Cells.Item(234, 6).Value  = 100114007
$ws.Cells.Item(234, 7).Value  = "Jengibre"
$ws.Cells.Item(234, 8).Value  = "Sin especificar"
$ws.Cells.Item(234, 9).Value  = "Primera"
$ws.Cells.Item(234, 10).Value = 220
$ws.Cells.Item(234, 11).Value = 20000
$ws.Cells.Item(234, 12).Value = 24000
$ws.Cells.Item(234, 13).Value = 22182
$ws.Cells.Item(234, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(234, 15).Value = "Perú"
$ws.Cells.Item(234, 16).Value = 1706
$ws.Cells.Item(234, 17).Value = 13
$ws.Cells.Item(234, 18).Value = "Hortaliza"
